$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the timestamp value stored in A14 (same day, slightly refined fractional time)
$ws.Range("A14").Value = 45868.54185008102

# Append new row 15 with the latest automated reading
$ws.Range("A15").Value = 45868.66690025864
$ws.Range("A15").Style = $ws.Range("A14").Style
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat

$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = 31
$ws.Range("D15").Value = 21.67
$ws.Range("E15").Value = 70.8
$ws.Range("F15").Value = 286.52
$ws.Range("G15").Value = 14.11
$ws.Range("H15").Value = "ESE"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "16:00:20"
